$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "36.180.88"
$ws.Range("E2").Value = "  -0.45%  "

# Row 3
$ws.Range("D3").Value = "2.005.33"
$ws.Range("E3").Value = "  -1.53%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
Set-TextValue $ws.Range("D5") "246.53"
$ws.Range("E5").Value = "  +0.85%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.644"
$ws.Range("E6").Value = "  -1.65%  "

# Row 7
Set-TextValue $ws.Range("D7") "61.43"
$ws.Range("E7").Value = "  +16.56%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
Set-TextValue $ws.Range("D9") "59.33"
$ws.Range("E9").Value = "  -2.35%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.369"
$ws.Range("E10").Value = "  +3.38%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0748"
$ws.Range("E11").Value = "  +1.98%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.105"
$ws.Range("E12").Value = "  -0.91%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.950"
$ws.Range("E13").Value = "  +3.76%  "

# Row 14
Set-TextValue $ws.Range("D14") "14.94"
$ws.Range("E14").Value = "  +4.20%  "

# Row 15
$ws.Range("D15").Value = "2.286.92"
$ws.Range("E15").Value = "  -2.11%  "

# Row 16
Set-TextValue $ws.Range("D16") "5.44"
$ws.Range("E16").Value = "  +2.28%  "

# Row 17
Set-TextValue $ws.Range("D17") "19.74"
$ws.Range("E17").Value = "  +18.40%  "

# Row 18
$ws.Range("D18").Value = "1.992.08"
$ws.Range("E18").Value = "  -2.51%  "

# Row 19
$ws.Range("D19").Value = "36.072.24"
$ws.Range("E19").Value = "  -0.56%  "

# Row 20
Set-TextValue $ws.Range("D20") "72.21"
$ws.Range("E20").Value = "  +2.07%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0858"
$ws.Range("E21").Value = "  +2.33%  "

# Row 22
$ws.Range("E22").Value = "  +3.25%  "

# Row 23
Set-TextValue $ws.Range("D23") "234.08"
$ws.Range("E23").Value = "  -0.31%  "

# Row 24
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D24") "2.64"
$ws.Range("E24").Value = "  +19.58%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D25") "1.00"
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("E26").Value = "  -2.97%  "

# Row 27
Set-TextValue $ws.Range("D27") "9.71"
$ws.Range("E27").Value = "  +7.81%  "

# Row 28
Set-TextValue $ws.Range("D28") "166.00"
$ws.Range("E28").Value = "  +1.89%  "

# Row 29
Set-TextValue $ws.Range("D29") "19.66"
$ws.Range("E29").Value = "  +0.20%  "

# Row 30
$ws.Range("E30").Value = "  +0.59%  "

# Row 31
Set-TextValue $ws.Range("D31") "5.11"
$ws.Range("E31").Value = "  +4.63%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.16"
$ws.Range("E32").Value = "  +0.01%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.0997"
$ws.Range("E33").Value = "  +15.76%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.0608"
$ws.Range("E34").Value = "  +4.33%  "

# Row 35
Set-TextValue $ws.Range("D35") "4.49"
$ws.Range("E35").Value = "  +3.78%  "

# Row 36
Set-TextValue $ws.Range("D36") "2.47"
$ws.Range("E36").Value = "  +14.17%  "

# Row 37
$ws.Range("E37").Value = "  -0.06%  "

# Row 38
$ws.Range("E38").Value = "  -1.57%  "

# Row 39
Set-TextValue $ws.Range("D39") "5.80"
$ws.Range("E39").Value = "  +18.61%  "

# Row 40
$ws.Range("E40").Value = "  +3.45%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.0961"
$ws.Range("E41").Value = "  +8.41%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D43") "0.0216"
$ws.Range("E43").Value = "  +2.75%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D44") "17.00"
$ws.Range("E44").Value = "  +10.48%  "

# Row 45
$ws.Range("E45").Value = "  +3.57%  "

# Row 46
Set-TextValue $ws.Range("D46") "94.42"
$ws.Range("E46").Value = "  +2.91%  "

# Row 47
Set-TextValue $ws.Range("D47") "7.82"
$ws.Range("E47").Value = "  +6.66%  "

# Row 48
$ws.Range("D48").Value = "1.369.72"
$ws.Range("E48").Value = "  +0.09%  "

# Row 49
Set-TextValue $ws.Range("D49") "2.92"
$ws.Range("E49").Value = "  +0.56%  "

# Row 50
$ws.Range("E50").Value = "  +5.77%  "

# Row 51
Set-TextValue $ws.Range("D51") "47.15"
$ws.Range("E51").Value = "  +6.84%  "
